$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the informational note in A1: append the extra sentence about "traces of".
$ws.Range("A1").Value = "Select also nuts when there is at least peanuts, almonds or hazelnuts. If whey, put it with milk. Also consider ""traces of""."

# Fill in allergen flags (0/1) for the first 10 labelled photos (rows 4-13).
# Columns: B=wheat, C=egg, D=milk, E=nuts, F=peanuts, G=almonds, H=hazelnuts,
#          I=soya, J=rye, K=oats, L=cheese, M=barley
$labels = @{
    4  = @(1,1,1,1,0,0,0,1,0,0,0,0)
    5  = @(1,1,1,1,0,0,0,1,0,0,0,0)
    6  = @(1,1,1,1,0,0,0,1,0,0,0,0)
    7  = @(1,1,1,1,0,0,0,1,0,0,0,0)
    8  = @(1,1,1,1,0,1,0,1,0,0,0,0)
    9  = @(0,0,1,1,0,1,1,1,0,0,0,0)
    10 = @(0,0,1,1,0,1,1,1,0,0,0,0)
    11 = @(0,0,1,1,0,1,1,1,0,0,0,0)
    12 = @(1,1,1,1,0,0,0,1,0,0,0,0)
    13 = @(1,1,1,1,0,0,1,0,0,0,0,0)
}

foreach ($row in 4..13) {
    $vals = $labels[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($row, 2 + $i).Value = $vals[$i]
    }
}

# Restore the view/selection to reflect the area just edited.
$ws.Range("N14").Select()
